$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Tnfsf10", "Tnfrsf11b", "FAPs", 3, 1, 17.54967, 52.64901, 0.9874120321709955, 0.9874120321709953, 3, 1, 2.095195666666667, 6.285587, 0.8546922300706357, 0.8546922300706358, 36.76999253543001, 330.92993281887, 0.8439333917748064, 0.8439333917748063),
    @("ECs",  "Tnfsf10", "Tnfrsf11b", "sCs",  3, 1, 17.54967, 52.64901, 0.9874120321709955, 0.9874120321709953, 2, 0.6666666666666666, 0.356208, 1.068624, 0.1453077699293643, 0.1453077699293643, 6.251332851360002, 56.26199566224, 0.143478640396189, 0.143478640396189),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.09171866666666667, 0.275156, 0.005160445469421788, 0.005160445469421787, 3, 1, 2.095195666666667, 6.285587, 0.8546922300706357, 0.8546922300706358, 0.1921685529524444, 1.729516976572, 0.004410592646418016, 0.004410592646418016),
    @("FAPs", "Tnfsf10", "Tnfrsf11b", "sCs",  1, 0.3333333333333333, 0.09171866666666667, 0.275156, 0.005160445469421788, 0.005160445469421787, 2, 0.6666666666666666, 0.356208, 1.068624, 0.1453077699293643, 0.1453077699293643, 0.032670922816, 0.294038305344, 0.0007498528230037713, 0.0007498528230037712),
    @("sCs",  "Tnfsf10", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.1320123333333333, 0.396037, 0.007427522359582914, 0.007427522359582914, 3, 1, 2.095195666666667, 6.285587, 0.8546922300706357, 0.8546922300706358, 0.2765916687465555, 2.489325018719, 0.006348245649411431, 0.006348245649411431),
    @("sCs",  "Tnfsf10", "Tnfrsf11b", "sCs",  1, 0.3333333333333333, 0.1320123333333333, 0.396037, 0.007427522359582914, 0.007427522359582914, 2, 0.6666666666666666, 0.356208, 1.068624, 0.1453077699293643, 0.1453077699293643, 0.047023849232, 0.423214643088, 0.001079276710171483, 0.001079276710171483)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
